# Update Value at Risk 95% (O), Portfolio Asset Weight (Q) and
# Portfolio Asset Allocation (R) figures on the Alternative, Bond and
# Equity sheets to reflect the new DIVIDEND_TYPE-driven calculations.

$wb = $excel.ActiveWorkbook

# --- Alternative sheet ---
$ws = $wb.Worksheets.Item("Alternative")
$ws.Range("O2").Value = -10.6
$ws.Range("Q2").Value = 4.5
$ws.Range("R2").Value = 450

$ws.Range("O3").Value = -21.69
$ws.Range("Q3").Value = 4.5
$ws.Range("R3").Value = 450

$ws.Range("O4").Value = -5.18
$ws.Range("Q4").Value = 21
$ws.Range("R4").Value = 2100

$ws.Range("O5").Value = -9.220000000000001

$ws.Range("O6").Value = -82.92

$ws.Range("O7").Value = -8.5

$ws.Range("O8").Value = -7.19

# --- Bond sheet ---
$ws = $wb.Worksheets.Item("Bond")
$ws.Range("O2").Value = -7.64
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0

$ws.Range("O3").Value = 3.82
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0

$ws.Range("O4").Value = -26.47
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0

$ws.Range("O5").Value = -16.96
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0

$ws.Range("O6").Value = -7.04
$ws.Range("Q6").Value = 12.22
$ws.Range("R6").Value = 1221.6

# --- Equity sheet ---
$ws = $wb.Worksheets.Item("Equity")
$ws.Range("O2").Value = -17.75
$ws.Range("Q2").Value = 9.24
$ws.Range("R2").Value = 923.62

$ws.Range("O3").Value = -15.72
$ws.Range("Q3").Value = 4.04
$ws.Range("R3").Value = 404.26

$ws.Range("O4").Value = -22.17
$ws.Range("Q4").Value = 5.59
$ws.Range("R4").Value = 558.54

$ws.Range("O5").Value = -16.75
$ws.Range("Q5").Value = 3.47
$ws.Range("R5").Value = 346.7

$ws.Range("O6").Value = -8.359999999999999
$ws.Range("Q6").Value = 4.82
$ws.Range("R6").Value = 481.57

$ws.Range("O7").Value = -14.01
$ws.Range("Q7").Value = 7.53
$ws.Range("R7").Value = 752.6900000000001

$ws.Range("O8").Value = -1.79
$ws.Range("Q8").Value = 23.11
$ws.Range("R8").Value = 2311.36

$ws.Range("O9").Value = -37.5

$wb.Save()
